$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): A1:O1 stay as-is; append 8 more header cells so
# the sheet now supports up to 5 ideas (idea3_* and idea4_* columns),
# matching the style ("s=1", bold/bordered/centered) already used by the
# other header cells.
$newHeaders = @(
    "idea3_ps_title", "idea3_ps_description", "idea3_title", "idea3_description",
    "idea4_ps_title", "idea4_ps_description", "idea4_title", "idea4_description"
)

$headerRange = $ws.Range("P1:W1")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, 16 + $i).Value = $newHeaders[$i]
}

# Copy the formatting of an existing header cell onto the new ones so they
# share the same style entry instead of creating a near-duplicate one.
$ws.Range("O1").Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2-4): the "team_id" values (T001/T002/T003) are dropped and
# every remaining value on each row shifts one column to the left. Rewrite
# each row explicitly with the new, shifted values; pad the rest of the row
# (up to column W) with empty cells.
$rowsData = @(
    @("Alpha Innovators", "PS-101", "Water Conservation", "Smart Water Saver", "Smart Drip System",
      "A low-cost IoT drip irrigation controller", "Air Quality", "Monitoring AQI in micro-climates",
      "AirSense", "Portable AQI monitor prototype"),
    @("Beta Builders", "PS-202", "Energy Efficiency", "Home Energy Optimizer", "HomeHub Energy Optimizer",
      "Central controller for home energy management", "Waste Management", "Smart sorting bin", "SmartSort",
      "Automated waste sorting using sensors", "Traffic Management", "Adaptive signals", "FlowSync",
      "Traffic flow adaptive algorithm"),
    @("Gamma Tech", "PS-303", "Healthcare Access", "TeleHealth Kiosk", "TeleHealth Kiosk for rural clinics",
      "Remote diagnostics and scheduling tool")
)

$lastCol = 23  # column W
for ($r = 0; $r -lt $rowsData.Length; $r++) {
    $excelRow = $r + 2
    $values = $rowsData[$r]
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($c -le $values.Length) {
            $ws.Cells.Item($excelRow, $c).Value = $values[$c - 1]
        } else {
            $ws.Cells.Item($excelRow, $c).Value = ""
        }
    }
}

$excel.CutCopyMode = $false

Write-Host "done"
